$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to keep text formatting so numeric-looking values
# (e.g. trailing zeros like "7.02", "2.00") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.869.60"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.948.22"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "553.61"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "133.39"
$ws.Range("E6").Value = "  +10.12%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +5.28%  "
$ws.Range("D9").Value = "2.943.90"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("D11").Value = "4.80"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  +4.43%  "
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").Value = "  +5.83%  "
$ws.Range("D14").Value = "32.80"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "3.432.57"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "6.96"
$ws.Range("E17").Value = "  +8.70%  "
$ws.Range("D18").Value = "2.944.20"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "57.800.50"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "416.11"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  +7.44%  "
$ws.Range("D23").Value = "13.41"
$ws.Range("E23").Value = "  +7.16%  "
$ws.Range("D24").Value = "7.02"
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("D25").Value = "79.21"
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "2.03"
$ws.Range("E29").Value = "  +6.94%  "
$ws.Range("D30").Value = "7.52"
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("D31").Value = "25.45"
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "0.0964"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.69"
$ws.Range("E34").Value = "  +6.57%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "0.946"
$ws.Range("E35").Value = "  +5.91%  "
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("D37").Value = "8.96"
$ws.Range("E37").Value = "  +7.40%  "
$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  +14.85%  "
$ws.Range("D39").Value = "48.25"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +16.72%  "
$ws.Range("D41").Value = "384.72"
$ws.Range("E41").Value = "  +8.06%  "
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").Value = "2.708.48"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "124.59"
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").Value = "22.91"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  +3.85%  "
